$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 20:46:11"
$ws1.Cells.Item(3, 1).Value = "Total filas: 364"
$ws1.Cells.Item(106, 1).Value = "10:05:51"
$ws1.Cells.Item(106, 2).Value = "11:52"
$ws1.Cells.Item(106, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(106, 4).Value = 107
$ws1.Cells.Item(106, 5).Value = "LP1912"
$ws1.Cells.Item(108, 1).Value = "11:47:17"
$ws1.Cells.Item(108, 2).Value = "11:52"
$ws1.Cells.Item(108, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(108, 4).Value = 5
$ws1.Cells.Item(108, 5).Value = "LP1912"
$ws1.Cells.Item(118, 1).Value = "11:34:59"
$ws1.Cells.Item(118, 2).Value = "12:09"
$ws1.Cells.Item(118, 3).Value = "15_ABASTO"
$ws1.Cells.Item(118, 4).Value = 35
$ws1.Cells.Item(118, 5).Value = "LP1912"
$ws1.Cells.Item(119, 1).Value = "11:34:59"
$ws1.Cells.Item(119, 2).Value = "12:09"
$ws1.Cells.Item(119, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(119, 4).Value = 35
$ws1.Cells.Item(119, 5).Value = "LP1912"
$ws1.Cells.Item(128, 1).Value = "11:47:17"
$ws1.Cells.Item(128, 2).Value = "12:32"
$ws1.Cells.Item(128, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(128, 4).Value = 45
$ws1.Cells.Item(128, 5).Value = "LP1912"
$ws1.Cells.Item(129, 1).Value = "10:37:52"
$ws1.Cells.Item(129, 2).Value = "12:32"
$ws1.Cells.Item(129, 3).Value = "14_ABASTO"
$ws1.Cells.Item(129, 4).Value = 115
$ws1.Cells.Item(129, 5).Value = "LP1912"
$ws1.Cells.Item(137, 1).Value = "11:47:17"
$ws1.Cells.Item(137, 2).Value = "12:37"
$ws1.Cells.Item(137, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(137, 4).Value = 50
$ws1.Cells.Item(137, 5).Value = "LP1912"
$ws1.Cells.Item(138, 1).Value = "11:52:01"
$ws1.Cells.Item(138, 2).Value = "12:37"
$ws1.Cells.Item(138, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(138, 4).Value = 45
$ws1.Cells.Item(138, 5).Value = "LP1912"
$ws1.Cells.Item(139, 1).Value = "11:34:59"
$ws1.Cells.Item(139, 2).Value = "12:47"
$ws1.Cells.Item(139, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(139, 4).Value = 73
$ws1.Cells.Item(139, 5).Value = "LP1912"
$ws1.Cells.Item(140, 1).Value = "11:34:59"
$ws1.Cells.Item(140, 2).Value = "12:47"
$ws1.Cells.Item(140, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(140, 4).Value = 73
$ws1.Cells.Item(140, 5).Value = "LP1912"
$ws1.Cells.Item(141, 1).Value = "11:34:59"
$ws1.Cells.Item(141, 2).Value = "12:47"
$ws1.Cells.Item(141, 3).Value = "14_ABASTO"
$ws1.Cells.Item(141, 4).Value = 73
$ws1.Cells.Item(141, 5).Value = "LP1912"
$ws1.Cells.Item(209, 1).Value = "14:53:07"
$ws1.Cells.Item(209, 2).Value = "15:53"
$ws1.Cells.Item(209, 3).Value = "10_OLMOS"
$ws1.Cells.Item(209, 4).Value = 60
$ws1.Cells.Item(209, 5).Value = "LP1912"
$ws1.Cells.Item(210, 1).Value = "13:56:11"
$ws1.Cells.Item(210, 2).Value = "15:53"
$ws1.Cells.Item(210, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(210, 4).Value = 117
$ws1.Cells.Item(210, 5).Value = "LP1912"
$ws1.Cells.Item(211, 1).Value = "13:56:11"
$ws1.Cells.Item(211, 2).Value = "15:53"
$ws1.Cells.Item(211, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(211, 4).Value = 117
$ws1.Cells.Item(211, 5).Value = "LP1912"
$ws1.Cells.Item(346, 1).Value = "19:55:23"
$ws1.Cells.Item(346, 2).Value = "20:55"
$ws1.Cells.Item(346, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(346, 4).Value = 60
$ws1.Cells.Item(346, 5).Value = "LP1912"
$ws1.Cells.Item(347, 1).Value = "19:48:11"
$ws1.Cells.Item(347, 2).Value = "20:55"
$ws1.Cells.Item(347, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(347, 4).Value = 67
$ws1.Cells.Item(347, 5).Value = "LP1912"
$ws1.Cells.Item(358, 1).Value = "19:55:23"
$ws1.Cells.Item(358, 2).Value = "21:34"
$ws1.Cells.Item(358, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(358, 4).Value = 99
$ws1.Cells.Item(358, 5).Value = "LP1912"
$ws1.Cells.Item(359, 1).Value = "19:48:11"
$ws1.Cells.Item(359, 2).Value = "21:34"
$ws1.Cells.Item(359, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(359, 4).Value = 106
$ws1.Cells.Item(359, 5).Value = "LP1912"
$ws1.Cells.Item(364, 1).Value = "20:46:10"
$ws1.Cells.Item(364, 2).Value = "21:50"
$ws1.Cells.Item(364, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(364, 4).Value = 64
$ws1.Cells.Item(364, 5).Value = "LP1912"
$ws1.Cells.Item(365, 1).Value = "20:32:13"
$ws1.Cells.Item(365, 2).Value = "22:04"
$ws1.Cells.Item(365, 3).Value = "15_ABASTO"
$ws1.Cells.Item(365, 4).Value = 92
$ws1.Cells.Item(365, 5).Value = "LP1912"
$ws1.Cells.Item(366, 1).Value = "20:32:13"
$ws1.Cells.Item(366, 2).Value = "22:11"
$ws1.Cells.Item(366, 3).Value = "14_ABASTO"
$ws1.Cells.Item(366, 4).Value = 99
$ws1.Cells.Item(366, 5).Value = "LP1912"
$ws1.Cells.Item(367, 1).Value = "20:46:10"
$ws1.Cells.Item(367, 2).Value = "22:34"
$ws1.Cells.Item(367, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(367, 4).Value = 108
$ws1.Cells.Item(367, 5).Value = "LP1912"
$ws1.Cells.Item(368, 1).Value = "20:46:10"
$ws1.Cells.Item(368, 2).Value = "22:34"
$ws1.Cells.Item(368, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(368, 4).Value = 108
$ws1.Cells.Item(368, 5).Value = "LP1912"
$ws1.Cells.Item(369, 1).Value = "20:46:10"
$ws1.Cells.Item(369, 2).Value = "22:44"
$ws1.Cells.Item(369, 3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(369, 4).Value = 118
$ws1.Cells.Item(369, 5).Value = "LP1912"

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 20:46:11"
$ws2.Cells.Item(3, 1).Value = "Total filas: 55"
$ws2.Cells.Item(47, 1).Value = "17:48:33"
$ws2.Cells.Item(47, 2).Value = "17:48"
$ws2.Cells.Item(47, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(47, 4).Value = 0
$ws2.Cells.Item(47, 5).Value = "LP1912"
$ws2.Cells.Item(48, 1).Value = "17:48:33"
$ws2.Cells.Item(48, 2).Value = "17:48"
$ws2.Cells.Item(48, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(48, 4).Value = 0
$ws2.Cells.Item(48, 5).Value = "LP1912"
$ws2.Cells.Item(60, 1).Value = "20:46:10"
$ws2.Cells.Item(60, 2).Value = "22:34"
$ws2.Cells.Item(60, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(60, 4).Value = 108
$ws2.Cells.Item(60, 5).Value = "LP1912"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 20:46:11"
$ws3.Cells.Item(3, 1).Value = "Total filas: 46"
$ws3.Cells.Item(49, 1).Value = "20:46:10"
$ws3.Cells.Item(49, 2).Value = "20:53"
$ws3.Cells.Item(49, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(49, 4).Value = 7
$ws3.Cells.Item(49, 5).Value = "L6203"
$ws3.Cells.Item(50, 1).Value = "20:12:03"
$ws3.Cells.Item(50, 2).Value = "21:30"
$ws3.Cells.Item(50, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(50, 4).Value = 78
$ws3.Cells.Item(50, 5).Value = "L6203"
$ws3.Cells.Item(51, 1).Value = "20:32:13"
$ws3.Cells.Item(51, 2).Value = "22:20"
$ws3.Cells.Item(51, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(51, 4).Value = 108
$ws3.Cells.Item(51, 5).Value = "L6173"
